$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Y2 value from "NA" to "Temperature-Rearing"
$ws.Range("Y2").Value = "Temperature-Rearing"

# Remove row 3 entirely (Nason Creek Lower 01 data), shifting dimension to A1:Y2
$ws.Rows.Item(3).Delete()
